$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Abraão"),
    @(2, "Rodrigo"),
    @(3, "Paulo"),
    @(4, "João"),
    @(5, "Jjj")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
